# Sprint 2 -> Sprint 3 planning update: fill in weekly progress tracking
# (SEM1-SEM4 = columns M:P) for the remaining backlog items, and move the
# active selection/scroll position to where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20's "PENDIENTE (Horas)" cell was a hard-coded number; restore it to
# the same shared formula used by the rest of the column (K - SUM(M:P)).
$ws.Range("L20").Formula = "=K20-(SUM(M20:P20))"

# Weekly hour allocations (SEM1=M, SEM2=N, SEM3=O, SEM4=P)
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = 2

$ws.Range("O10").Value = 2
$ws.Range("P10").Value = 2

$ws.Range("O11").Value = 2
$ws.Range("P11").Value = 2

$ws.Range("P15").Value = 4

$ws.Range("O17").Value = 2
$ws.Range("P17").Value = 2

$ws.Range("N18").Value = 2
$ws.Range("O18").Value = 2

$ws.Range("N19").Value = 2
$ws.Range("O19").Value = 2

$ws.Range("M20").Value = 2
$ws.Range("N20").Value = 2
$ws.Range("O20").Value = 2

$ws.Range("P21").Value = 2

$ws.Range("M22").Value = 16

$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 2

$ws.Range("N24").Value = 2
$ws.Range("O24").Value = 2

$ws.Range("O25").Value = 2
$ws.Range("P25").Value = 2

$ws.Range("P26").Value = 2

$ws.Range("P27").Value = 2

$ws.Range("P28").Value = 2

$ws.Range("M29").Value = 8
$ws.Range("P29").Value = 4

$ws.Range("P30").Value = 2

$ws.Range("P31").Value = 4

# Scroll the view down a row and park the selection on the last edited cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("P31").Select()
